# Logic For Screenplay Reader.xlsx - "GOt the generic function working I think."
#
# - Adds three new annotation cells to the existing sheet (SCENE / CHARACTER / PARENTH)
# - Renames the original sheet to "LOGIC of First Attempt"
# - Adds a new "Sheet2" (placed after the first sheet, becomes the active tab) that
#   lists out the token labels used by the parser logic.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- annotate the existing logic sheet -----------------------------------
$ws1.Range("A10").Value = "SCENE"
$ws1.Range("A17").Value = "CHARACTER"
$ws1.Range("A28").Value = "PARENTH"

# --- create the second sheet ----------------------------------------------
# Copy+clear (rather than Worksheets.Add) so the new sheet inherits the same
# worksheet-level formatting scaffolding as the original.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Clear()
$ws2.Name = "Sheet2"

$ws1.Name = "LOGIC of First Attempt"

$ws2.Range("A2").Value = "BLANK_LINE"
$ws2.Range("A3").Value = "CUT"
$ws2.Range("A4").Value = "CAMERA DIRECTION "
$ws2.Range("A5").Value = "TRANSISTION "
$ws2.Range("A7").Value = "PAGE_NUMBER"
$ws2.Range("A6").Value = "CONTINUED_PAGE "
$ws2.Range("A1").Value = "Determine with REgWhy"
$ws2.Range("A8").Value = "SCENE"
$ws2.Range("A9").Value = "CHARACTER"

# --- selections / active tab ----------------------------------------------
# Leave sheet1's selection where it was left (C7), then select Sheet2's A10
# last so Sheet2 ends up the active/visible tab.
$ws1.Range("C7").Select() | Out-Null
$ws2.Range("A10").Select() | Out-Null
